# Edit script: applies the changes described by the target diff.
#
# 1) Three tables (on slides 14, 15 and 16) have their table style swapped
#    from the custom "Table_0" style ({225F68CB-473D-46D8-BD90-08B3FA5F88A7})
#    to the built-in "Medium Style 2 - Accent 1" style
#    ({6F7A0A0E-0110-488E-9D02-78BF985CBE24}).
#
# 2) The deck's theme (ppt/theme/theme2.xml, the one actually applied to the
#    slide master / slides) is switched from the "Integral" / "Red Violet"
#    colour scheme back to the stock "Office" colour scheme (i.e. the colour
#    values that otherwise live in ppt/theme/theme1.xml, which is only used
#    by the notes master).

$p = $ppt.ActivePresentation

# --- 1) Table styles -------------------------------------------------
$newStyleId = "{6F7A0A0E-0110-488E-9D02-78BF985CBE24}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2) Theme colour scheme -------------------------------------------
# Office colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink),
# in MsoThemeColorSchemeIndex order (1..12).
$officeColors = @(
    0x000000,
    0xFFFFFF,
    0x44546A,
    0xE7E6E6,
    0x5B9BD5,
    0xED7D31,
    0xA5A5A5,
    0xFFC000,
    0x4472C4,
    0x70AD47,
    0x0563C1,
    0x954F72
)

$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $bgr = $officeColors[$i - 1]
    $r = ($bgr -band 0xFF0000) -shr 16
    $g = ($bgr -band 0x00FF00) -shr 8
    $b = ($bgr -band 0x0000FF)
    $colorScheme.Item($i).RGB = $r + ($g * 256) + ($b * 65536)
}
